$wb = $excel.ActiveWorkbook

# Add the new "Kosten" worksheet after the last existing sheet
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Kosten"

# Header row
$ws.Cells.Item(1, 1).Value = "Spieler"
$ws.Cells.Item(1, 2).Value = "Teilnahmen"
$ws.Cells.Item(1, 3).Value = "Minuten"
$ws.Cells.Item(1, 4).Value = "Summe"

# Data rows
$data = @(
    @("Frank Petermann", 18, 1950, 174.97),
    @("Dirk Kistner", 15, 900, 131.25),
    @("Thomas Bretschneider", 10, 900, 131.19999999999999),
    @("MARTIN LANGE", 14, 1680, 122.5),
    @("Liam Wilde", 10, 990, 118.09),
    @("Björn Junker", 13, 1560, 113.75),
    @("lars staubermann", 13, 1560, 113.75),
    @("Manfred Grell", 11, 900, 111.53),
    @("Michael Rabehl", 10, 1050, 109.35),
    @("Oliver Böss", 9, 750, 109.34),
    @("Frank Koller", 10, 780, 107.16),
    @("Sebastian Braune", 10, 1080, 104.98),
    @("Thomas Grüneberg", 12, 840, 102.8),
    @("Carsten Gambal", 12, 1230, 100.61),
    @("Matthias Duddek", 11, 660, 96.25),
    @("Jan Pappenheim", 11, 1050, 94.039999999999992),
    @("Bernd Robioneck", 10, 900, 85.28),
    @("Andreas Dank", 11, 1080, 83.11),
    @("Peter Plähn", 9, 1080, 78.75),
    @("Anke Ihde", 10, 1080, 78.739999999999995),
    @("Markus Muench", 10, 810, 72.17),
    @("Jens Krause", 11, 990, 72.16),
    @("Jörg Peters", 10, 840, 69.97999999999999),
    @("Michael Bock", 10, 960, 69.97999999999999),
    @("Bernd Sotzek", 9, 750, 67.8),
    @("Ralf Colditz", 9, 690, 67.8),
    @("Jens Hafner", 9, 840, 65.61),
    @("Jürgen Hansen", 10, 900, 65.599999999999994),
    @("Karsten Usinger", 8, 720, 65.599999999999994),
    @("Lena Meiß", 7, 840, 61.25),
    @("mohamad albadry", 7, 840, 61.25),
    @("Martina Schmidt", 8, 840, 61.24),
    @("Heiko Thomsen", 7, 450, 59.06),
    @("Lorenz Kramp", 7, 570, 54.68),
    @("arndt stueber", 6, 360, 52.5),
    @("kai schroeder", 6, 360, 52.5),
    @("Tobias Kahl", 6, 600, 48.12),
    @("Torsten Bartel", 4, 390, 48.11),
    @("Patrick Buehrsch", 5, 300, 43.75),
    @("patrick buehrsch", 4, 240, 35),
    @("Gunnar Brix", 5, 480, 34.989999999999988),
    @("Kai Schröder", 3, 180, 26.25),
    @("Kerstin Baarck", 3, 360, 26.25),
    @("liam wilde", 2, 180, 26.24),
    @("wolfgang aleksik", 4, 360, 26.24),
    @("torsten bartel", 2, 180, 19.68),
    @("Mohamad Albadry", 2, 240, 17.5),
    @("tobias kahl", 2, 180, 17.5),
    @("lorenz kramp", 2, 150, 15.31),
    @("Arndt Stüber", 2, 180, 13.12),
    @("Wolfgang Aleksik", 2, 180, 13.12),
    @("karsten usinger", 1, 90, 13.12),
    @("heiko thomsen", 1, 60, 8.75),
    @("gunnar brix", 1, 90, 6.56),
    @("jens hafner", 1, 90, 6.56),
    @("martina schmidt", 1, 90, 6.56),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Restore the originally active sheet ("Herren 40-50-60")
$lastSheet.Activate()
